# Latest calibration for Fiji
# Updates the calibrated constant values on the "constants" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value  = 10.03626685000244
$ws.Range("B3").Value  = 0.09999999999999998
$ws.Range("B4").Value  = 1845.010175260038
$ws.Range("B5").Value  = 18409.4206532835
$ws.Range("B8").Value  = 0.624659011346387
$ws.Range("B9").Value  = 0.7262650305460934
$ws.Range("B10").Value = 2.816856563164909
$ws.Range("B11").Value = 0.2902694698354936
